$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test description text in B2
$ws.Range("B2").Value = "Check Edit Lead Info Button Enabled from Lead Info tab (User is creator of the nominated GF to DB Lead)"

# Widen column B and let Excel manage the fit state (no longer best-fit, explicit width)
$ws.Columns.Item(2).ColumnWidth = 97.5

# Update the sheet view: clear the frozen/scrolled top-left cell and move selection to column C (whole column)
$ws.Range("C1:C1048576").Select()
